# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions".
#
# Most touched cells (links, coin names, percentage strings) are plain
# text already, so a normal .Value assignment round-trips them exactly.
# A handful of Price cells look like bare numbers ("1.00", "0.165", ...)
# and Excel's COM layer would silently coerce those to numeric cells,
# rounding/retyping them. For those we enter the value with a leading
# apostrophe (forces text entry, like typing '0.165 into a cell) and then
# reset .Style so the quote-prefix flag doesn't leave a stray style on
# the cell - matching the workbook's original unstyled text cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.411.67'
$ws.Range("E2").Value = '  -1.34%  '
$ws.Range("D3").Value = '3.849.11'
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("E4").Value = '  -0.12%  '
$c = $ws.Range("D5")
$c.Value = '''601.94'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = '''168.97'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").Value = '3.848.65'
$ws.Range("E7").Value = '  -1.07%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -0.93%  '
$c = $ws.Range("D10")
$c.Value = '''0.165'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.95%  '
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("E12").Value = '  -2.05%  '
$c = $ws.Range("D13")
$c.Value = '''0.0000266'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +4.63%  '
$c = $ws.Range("D14")
$c.Value = '''37.12'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -2.73%  '
$ws.Range("D15").Value = '4.494.58'
$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("D16").Value = '3.846.80'
$ws.Range("E16").Value = '  -1.12%  '
$ws.Range("D17").Value = '68.504.05'
$ws.Range("E17").Value = '  -1.34%  '
$c = $ws.Range("D18")
$c.Value = '''18.53'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("E19").Value = '  -3.05%  '
$ws.Range("E20").Value = '  -1.05%  '
$c = $ws.Range("D21")
$c.Value = '''11.21'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.68%  '
$c = $ws.Range("D22")
$c.Value = '''470.81'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -3.80%  '
$c = $ws.Range("D23")
$c.Value = '''0.733'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.13%  '
$ws.Range("E24").Value = '  -3.27%  '
$c = $ws.Range("D25")
$c.Value = '''83.49'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("E26").Value = '  -2.28%  '
$c = $ws.Range("D27")
$c.Value = '''12.12'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.86%  '
$c = $ws.Range("D28")
$c.Value = '''10.23'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +1.29%  '
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("E30").Value = '  -0.82%  '
$ws.Range("D31").Value = '3.998.93'
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("E32").Value = '  -1.03%  '
$ws.Range("E33").Value = '  -1.46%  '
$c = $ws.Range("D34")
$c.Value = '''2.30'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -3.24%  '
$ws.Range("E35").Value = '  -2.81%  '
$ws.Range("D36").Value = '3.815.05'
$ws.Range("E36").Value = '  -1.05%  '
$c = $ws.Range("D37")
$c.Value = '''3.79'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +12.36%  '
$ws.Range("E38").Value = '  -1.93%  '
$ws.Range("E39").Value = '  -1.02%  '
$ws.Range("E40").Value = '  -2.01%  '
$c = $ws.Range("D41")
$c.Value = '''5.94'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -2.37%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("E43").Value = '  -3.09%  '
$ws.Range("E44").Value = '  -3.91%  '
$c = $ws.Range("D45")
$c.Value = '''8.71'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.44%  '
$c = $ws.Range("D46")
$c.Value = '''417.79'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -3.94%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Range("D47")
$c.Value = '''1.00'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$c = $ws.Range("D48")
$c.Value = '''0.000292'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +6.05%  '
$ws.Range("E49").Value = '  -2.19%  '
$c = $ws.Range("D50")
$c.Value = '''142.02'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.55%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D51")
$c.Value = '''26.09'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +4.07%  '
